# Update 北京-漫展信息.xlsx: refresh "想去人数" (F column) counts across all
# sheets, and swap one cover image URL (I32 on 展览, I36 on 全部类型), to
# match the regenerated gh-pages output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7652
$ws1.Range("F3").Value  = 7652
$ws1.Range("F5").Value  = 7835
$ws1.Range("F8").Value  = 29
$ws1.Range("F9").Value  = 6572
$ws1.Range("F10").Value = 3354
$ws1.Range("F12").Value = 3708
$ws1.Range("F14").Value = 43
$ws1.Range("F16").Value = 62
$ws1.Range("F17").Value = 55
$ws1.Range("F18").Value = 463
$ws1.Range("F19").Value = 8
$ws1.Range("F20").Value = 15
$ws1.Range("F21").Value = 308
$ws1.Range("F22").Value = 323
$ws1.Range("F23").Value = 3812
$ws1.Range("F26").Value = 953
$ws1.Range("F28").Value = 1460
$ws1.Range("F29").Value = 79
$ws1.Range("F30").Value = 52
$ws1.Range("F31").Value = 2732
$ws1.Range("F32").Value = 1780
$ws1.Range("F35").Value = 56
$ws1.Range("F36").Value = 3625
$ws1.Range("F37").Value = 298
$ws1.Range("F38").Value = 276
$ws1.Range("F41").Value = 532
$ws1.Range("F42").Value = 1404
$ws1.Range("F43").Value = 242
$ws1.Range("F44").Value = 548
$ws1.Range("F45").Value = 634
$ws1.Range("I32").Value = "//i2.hdslb.com/bfs/openplatform/202405/9CAdQvG71716812495452.jpeg"

# --- Sheet: 演出 (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 239
$ws2.Range("F13").Value = 88
$ws2.Range("F17").Value = 14

# --- Sheet: 本地生活 (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 131

# --- Sheet: 全部类型 (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 131
$ws4.Range("F5").Value  = 7652
$ws4.Range("F6").Value  = 7652
$ws4.Range("F8").Value  = 7835
$ws4.Range("F10").Value = 29
$ws4.Range("F11").Value = 6572
$ws4.Range("F12").Value = 3354
$ws4.Range("F14").Value = 3708
$ws4.Range("F16").Value = 43
$ws4.Range("F18").Value = 62
$ws4.Range("F19").Value = 55
$ws4.Range("F20").Value = 463
$ws4.Range("F21").Value = 308
$ws4.Range("F23").Value = 323
$ws4.Range("F24").Value = 3812
$ws4.Range("F30").Value = 953
$ws4.Range("F32").Value = 1460
$ws4.Range("F33").Value = 79
$ws4.Range("F34").Value = 52
$ws4.Range("F35").Value = 2732
$ws4.Range("F36").Value = 1780
$ws4.Range("F39").Value = 88
$ws4.Range("F40").Value = 3625
$ws4.Range("F41").Value = 298
$ws4.Range("F42").Value = 276
$ws4.Range("F45").Value = 532
$ws4.Range("F46").Value = 1404
$ws4.Range("F47").Value = 242
$ws4.Range("F49").Value = 548
$ws4.Range("F50").Value = 634
$ws4.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202405/9CAdQvG71716812495452.jpeg"
